$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12..77 down to 13..78
# (Excel's Insert() copies the formatting of the row above, which already
# gives D12 the correct date number-format style.)
$ws.Rows("12:12").Insert()

# Populate the new record in row 12 with the same constant descriptor
# columns used throughout the sheet, plus the new observation's values.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 44547
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100101
$ws.Range("H12").Value = "Berries"
$ws.Range("I12").Value = 100101001
$ws.Range("J12").Value = "Arándano (blue)"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 500
$ws.Range("N12").Value = 2600
$ws.Range("O12").Value = 2600
$ws.Range("P12").Value = 2600
$ws.Range("Q12").Value = "$/kilo"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 2600
$ws.Range("T12").Value = 1
